$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1424.6666
$ws.Range("H31").Value = 1131.0834
$ws.Range("I31").Value = 521.625
$ws.Range("J31").Value = 2350
$ws.Range("K31").Value = 1564.875
$ws.Range("L31").Value = 7050
$ws.Range("M31").Value = -1334.875
$ws.Range("N31").Value = -7510
$ws.Range("H38").Value = 73.888885
$ws.Range("I38").Value = 73.888885
$ws.Range("K38").Value = 221.666655
$ws.Range("M38").Value = 150.333345
$ws.Range("H58").Value = 1001171.5
$ws.Range("I58").Value = 251.66667
$ws.Range("J58").Value = 1430137.1
$ws.Range("K58").Value = 755.00001
$ws.Range("L58").Value = 4290411.300000001
$ws.Range("M58").Value = -605.00001
$ws.Range("N58").Value = -4290711.300000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3803.1428
$ws.Range("I2").Value = 3803.1428
$ws.Range("K2").Value = 3803.1428
$ws.Range("M2").Value = -3690.1428
$ws.Range("H32").Value = 2251.82
$ws.Range("I32").Value = 2166.1428
$ws.Range("J32").Value = 6450
$ws.Range("K32").Value = 2166.1428
$ws.Range("L32").Value = 6450
$ws.Range("M32").Value = -1879.1428
$ws.Range("N32").Value = -7024
$ws.Range("H33").Value = 12000
$ws.Range("J33").Value = 12000
$ws.Range("L33").Value = 12000
$ws.Range("N33").Value = -12658
$ws.Range("H45").Value = 1248.6
$ws.Range("I45").Value = 1023.73334
$ws.Range("J45").Value = 1585.9
$ws.Range("K45").Value = 1023.73334
$ws.Range("L45").Value = 1585.9
$ws.Range("M45").Value = -646.73334
$ws.Range("N45").Value = -2339.9
$ws.Range("H61").Value = 709.94116
$ws.Range("I61").Value = 613.6774
$ws.Range("J61").Value = 1704.6666
$ws.Range("K61").Value = 613.6774
$ws.Range("L61").Value = 1704.6666
$ws.Range("M61").Value = -401.6774
$ws.Range("N61").Value = -2128.6666
$ws.Range("H97").Value = 1075.9744
$ws.Range("I97").Value = 751.06665
$ws.Range("K97").Value = 751.06665
$ws.Range("M97").Value = -255.06665
$ws.Range("H102").Value = 1575
$ws.Range("I102").Value = 1100
$ws.Range("K102").Value = 1100
$ws.Range("M102").Value = 522
$ws.Range("H116").Value = 3803.1428
$ws.Range("I116").Value = 3803.1428
$ws.Range("K116").Value = 3803.1428
$ws.Range("M116").Value = -1509.1428
$ws.Range("H136").Value = 709.94116
$ws.Range("I136").Value = 613.6774
$ws.Range("J136").Value = 1704.6666
$ws.Range("K136").Value = 1841.0322
$ws.Range("L136").Value = 5113.9998
$ws.Range("M136").Value = 708.9677999999999
$ws.Range("N136").Value = -10213.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3803.1428
$ws.Range("I3").Value = 3803.1428
$ws.Range("K3").Value = 3803.1428
$ws.Range("M3").Value = -3689.1428
$ws.Range("H32").Value = 49000
$ws.Range("J32").Value = 49000
$ws.Range("L32").Value = 49000
$ws.Range("N32").Value = -49768
$ws.Range("H99").Value = 1183.1
$ws.Range("I99").Value = 1136.6666
$ws.Range("J99").Value = 1252.75
$ws.Range("K99").Value = 1136.6666
$ws.Range("L99").Value = 1252.75
$ws.Range("M99").Value = 361.3334
$ws.Range("N99").Value = -4248.75
$ws.Range("H134").Value = 1581.5588
$ws.Range("I134").Value = 1003.7692
$ws.Range("J134").Value = 3459.375
$ws.Range("K134").Value = 3011.3076
$ws.Range("L134").Value = 10378.125
$ws.Range("M134").Value = -476.3076000000001
$ws.Range("N134").Value = -15448.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H138").Value = 64957
$ws.Range("J138").Value = 64957
$ws.Range("L138").Value = 64957
$ws.Range("N138").Value = -75237

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 2500
$ws.Range("J110").Value = 4000
$ws.Range("L110").Value = 12000
$ws.Range("N110").Value = -20180
$ws.Range("H117").Value = 1214.5
$ws.Range("I117").Value = 1214.5
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 3643.5
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = -201.5
$ws.Range("N117").ClearContents()
$ws.Range("H119").Value = 4578.3887
$ws.Range("I119").Value = 888.1667
$ws.Range("J119").Value = 6423.5
$ws.Range("K119").Value = 2664.5001
$ws.Range("L119").Value = 19270.5
$ws.Range("M119").Value = 2173.4999
$ws.Range("N119").Value = -28946.5
$ws.Range("H120").Value = 2605.625
$ws.Range("I120").Value = 2605.625
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 7816.875
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -2978.875
$ws.Range("N120").ClearContents()
$ws.Range("H131").Value = 3885.4055
$ws.Range("J131").Value = 4656
$ws.Range("L131").Value = 13968
$ws.Range("N131").Value = -24048

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1851.0769
$ws.Range("I132").Value = 1499.8276
$ws.Range("J132").Value = 2869.7
$ws.Range("K132").Value = 4499.4828
$ws.Range("L132").Value = 8609.099999999999
$ws.Range("M132").Value = -1969.4828
$ws.Range("N132").Value = -13669.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("H61").Value = 2002.5
$ws.Range("J61").Value = 2005
$ws.Range("L61").Value = 2005
$ws.Range("N61").Value = -2409
$ws.Range("H93").Value = 7460.875
$ws.Range("I93").Value = 10091.363
$ws.Range("K93").Value = 10091.363
$ws.Range("M93").Value = -8843.362999999999
$ws.Range("H113").Value = 2002.5
$ws.Range("J113").Value = 2005
$ws.Range("L113").Value = 2005
$ws.Range("N113").Value = -6345
$ws.Range("H132").Value = 4940.033
$ws.Range("I132").Value = 4976
$ws.Range("J132").Value = 4870.4194
$ws.Range("K132").Value = 14928
$ws.Range("L132").Value = 14611.2582
$ws.Range("M132").Value = -12398
$ws.Range("N132").Value = -19671.2582
$ws.Range("H136").Value = 1979.909
$ws.Range("I136").Value = 1978.7142
$ws.Range("J136").Value = 2005
$ws.Range("K136").Value = 5936.142599999999
$ws.Range("L136").Value = 2005
$ws.Range("M136").Value = -3386.142599999999
$ws.Range("N136").Value = -11115
